$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.636.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.564.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.82"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.98%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0897"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.788.20"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.567.98"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.674.94"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.51"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.41"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.31"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0681"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.92"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.05"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.46"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.107"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.08%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.24"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0457"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.91%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.18"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.403.93"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.68"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.94"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.516"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.89"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.700.04"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -6.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.76"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.59"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.09%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.41%  "
